# Update database and change read_price algorithm.
#
# The yearly cost-report sheet keeps a rolling 5-period window (columns
# E:I) for every metric row. A new period ("1401/12") was appended, so
# every row's window slides one column to the left (E<-F, F<-G, G<-H,
# H<-I) and the freed-up column I receives the new period's figure
# (or, for the two header rows per block, the new period label).
#
# A handful of "amount" rows (E.g. row 29, 36, 43, ...) used to show a
# literal "-" placeholder in column E for periods that predated the
# read_price algorithm; now that the algorithm back-fills real figures
# for every period, those placeholders disappear along with the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newYearLabel = "دوازده ماهه منتهی به 1401/12"

# Rows whose E:I block is a year-header (label row) rather than data.
$headerRows = @(8,27,34,41,48,55,62,69,76,83,89,95,101,107)

# Every other row that carries data across E:I, mapped to the new
# (sixth-period) value that lands in column I after the shift.
$newIValue = @{
    10  = 13820299
    11  = 136347
    12  = 1847234
    13  = 15803880
    14  = 0
    15  = 15803880
    16  = 0
    17  = 0
    18  = 15803880
    19  = 1258329
    20  = -788330
    21  = 16273879
    22  = 0
    23  = 16273879
    29  = 3473799
    30  = 3473799
    36  = 103091984
    37  = 103091984
    43  = 97127177
    44  = 97127177
    50  = 9438606
    51  = 9438606
    57  = 507697
    58  = 507697
    64  = 14475963
    65  = 14475963
    71  = 13820299
    72  = 13820299
    78  = 1163361
    79  = 1163361
    85  = 146150
    91  = 140418
    97  = 142291
    103 = 123256
    109 = 0
    110 = 0
    111 = 0
    112 = 0
    113 = 0
    114 = 150814
    115 = 256579
    116 = 615873
    117 = 0
    118 = 823968
    119 = 1847234
}

foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 6).Value2   # E <- F
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 7).Value2   # F <- G
    $ws.Cells.Item($r, 7).Value2 = $ws.Cells.Item($r, 8).Value2   # G <- H
    $ws.Cells.Item($r, 8).Value2 = $ws.Cells.Item($r, 9).Value2   # H <- I
    $ws.Cells.Item($r, 9).Value2 = $newYearLabel                  # I <- new period
}

foreach ($r in $newIValue.Keys) {
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 6).Value2   # E <- F
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 7).Value2   # F <- G
    $ws.Cells.Item($r, 7).Value2 = $ws.Cells.Item($r, 8).Value2   # G <- H
    $ws.Cells.Item($r, 8).Value2 = $ws.Cells.Item($r, 9).Value2   # H <- I
    $ws.Cells.Item($r, 9).Value2 = $newIValue[$r]                 # I <- new period
}
